$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8, shifting existing rows 8-25 down to 9-26.
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the latest weekly price entry.
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value = 44414
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100108
$ws.Cells.Item(8, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(8, 9).Value = 100108001
$ws.Cells.Item(8, 10).Value = "Guayaba"
$ws.Cells.Item(8, 11).Value = "Sin especificar"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 160
$ws.Cells.Item(8, 14).Value = 1300
$ws.Cells.Item(8, 15).Value = 1400
$ws.Cells.Item(8, 16).Value = 1350
$ws.Cells.Item(8, 17).Value = "$/kilo (en caja de 10 kilos )"
$ws.Cells.Item(8, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 19).Value = 1350
$ws.Cells.Item(8, 20).Value = 1
